$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 336, 337
$ws.Range("B336").Value = 4591717
$ws.Range("F336").Value = "Orgryte IS"
$ws.Range("G336").Value = "Orebro SK"
$ws.Range("H336").Value = 2
$ws.Range("I336").Value = 0
$ws.Range("J336").Value = "H"
$ws.Range("K336").Value = 2.55
$ws.Range("L336").Value = 3.5
$ws.Range("M336").Value = 2.55
$ws.Range("N336").Value = 2
$ws.Range("O336").Value = 3.75
$ws.Range("P336").Value = 3.4
$ws.Range("Q336").Value = -0.5
$ws.Range("R336").Value = 2.025
$ws.Range("S336").Value = 1.825
$ws.Range("T336").Value = 2.5
$ws.Range("U336").Value = 2
$ws.Range("V336").Value = 1.85
$ws.Range("W336").Value = 1
$ws.Range("X336").Value = -1
$ws.Range("Y336").Value = -1
$ws.Range("Z336").Value = 1.025
$ws.Range("AA336").Value = -1
$ws.Range("AB336").Value = -1
$ws.Range("AC336").Value = 0.8500000000000001
$ws.Range("B337").Value = 4587838
$ws.Range("F337").Value = "Halmstad"
$ws.Range("G337").Value = "Jonkopings Sodra"
$ws.Range("H337").Value = 3
$ws.Range("I337").Value = 0
$ws.Range("J337").Value = "H"
$ws.Range("K337").Value = 1.533
$ws.Range("L337").Value = 4.2
$ws.Range("M337").Value = 5.75
$ws.Range("N337").Value = 1.4
$ws.Range("O337").Value = 4.75
$ws.Range("P337").Value = 7
$ws.Range("Q337").Value = -1.25
$ws.Range("R337").Value = 2
$ws.Range("S337").Value = 1.85
$ws.Range("T337").Value = 2.75
$ws.Range("U337").Value = 1.975
$ws.Range("V337").Value = 1.875
$ws.Range("W337").Value = 0.3999999999999999
$ws.Range("X337").Value = -1
$ws.Range("Y337").Value = -1
$ws.Range("Z337").Value = 1
$ws.Range("AA337").Value = -1
$ws.Range("AB337").Value = 0.4875
$ws.Range("AC337").Value = -0.5

# Rows 344, 345
$ws.Range("B344").Value = 4587840
$ws.Range("F344").Value = "Landskrona BoIS"
$ws.Range("G344").Value = "Halmstad"
$ws.Range("H344").Value = 2
$ws.Range("I344").Value = 1
$ws.Range("J344").Value = "H"
$ws.Range("K344").Value = 3.2
$ws.Range("L344").Value = 3.4
$ws.Range("M344").Value = 2
$ws.Range("N344").Value = 4
$ws.Range("O344").Value = 3.25
$ws.Range("P344").Value = 1.833
$ws.Range("Q344").Value = 0.5
$ws.Range("R344").Value = 1.95
$ws.Range("S344").Value = 1.9
$ws.Range("T344").Value = 2.25
$ws.Range("U344").Value = 1.85
$ws.Range("V344").Value = 2
$ws.Range("W344").Value = 3
$ws.Range("X344").Value = -1
$ws.Range("Y344").Value = -1
$ws.Range("Z344").Value = 0.95
$ws.Range("AA344").Value = -1
$ws.Range("AB344").Value = 0.8500000000000001
$ws.Range("AC344").Value = -1
$ws.Range("B345").Value = 4591716
$ws.Range("F345").Value = "Orebro SK"
$ws.Range("G345").Value = "Trelleborgs FF"
$ws.Range("H345").Value = 2
$ws.Range("I345").Value = 2
$ws.Range("J345").Value = "D"
$ws.Range("K345").Value = 2.05
$ws.Range("L345").Value = 3.4
$ws.Range("M345").Value = 3.3
$ws.Range("N345").Value = 2.55
$ws.Range("O345").Value = 3.3
$ws.Range("P345").Value = 2.6
$ws.Range("Q345").Value = 0
$ws.Range("R345").Value = 1.95
$ws.Range("S345").Value = 1.9
$ws.Range("T345").Value = 2.5
$ws.Range("U345").Value = 1.85
$ws.Range("V345").Value = 2
$ws.Range("W345").Value = -1
$ws.Range("X345").Value = 2.3
$ws.Range("Y345").Value = -1
$ws.Range("Z345").Value = 0
$ws.Range("AA345").Value = -0
$ws.Range("AB345").Value = 0.8500000000000001
$ws.Range("AC345").Value = -1

# Rows 352, 353
$ws.Range("B352").Value = 4587843
$ws.Range("F352").Value = "Orgryte IS"
$ws.Range("G352").Value = "Brommapojkarna"
$ws.Range("H352").Value = 1
$ws.Range("I352").Value = 1
$ws.Range("J352").Value = "D"
$ws.Range("K352").Value = 2.625
$ws.Range("L352").Value = 3.4
$ws.Range("M352").Value = 2.375
$ws.Range("N352").Value = 3
$ws.Range("O352").Value = 3.6
$ws.Range("P352").Value = 1.95
$ws.Range("Q352").Value = 0.5
$ws.Range("R352").Value = 1.8
$ws.Range("S352").Value = 2.05
$ws.Range("T352").Value = 3
$ws.Range("U352").Value = 2.025
$ws.Range("V352").Value = 1.825
$ws.Range("W352").Value = -1
$ws.Range("X352").Value = 2.6
$ws.Range("Y352").Value = -1
$ws.Range("Z352").Value = 0.8
$ws.Range("AA352").Value = -1
$ws.Range("AB352").Value = -1
$ws.Range("AC352").Value = 0.825
$ws.Range("B353").Value = 4587976
$ws.Range("F353").Value = "Norrby IF"
$ws.Range("G353").Value = "Jonkopings Sodra"
$ws.Range("H353").Value = 1
$ws.Range("I353").Value = 1
$ws.Range("J353").Value = "D"
$ws.Range("K353").Value = 1.75
$ws.Range("L353").Value = 3.4
$ws.Range("M353").Value = 4.2
$ws.Range("N353").Value = 1.6
$ws.Range("O353").Value = 3.6
$ws.Range("P353").Value = 5.25
$ws.Range("Q353").Value = -0.75
$ws.Range("R353").Value = 1.825
$ws.Range("S353").Value = 2.025
$ws.Range("T353").Value = 2.5
$ws.Range("U353").Value = 1.975
$ws.Range("V353").Value = 1.875
$ws.Range("W353").Value = -1
$ws.Range("X353").Value = 2.6
$ws.Range("Y353").Value = -1
$ws.Range("Z353").Value = -1
$ws.Range("AA353").Value = 1.025
$ws.Range("AB353").Value = -1
$ws.Range("AC353").Value = 0.875

# Rows 368, 369
$ws.Range("B368").Value = 4587854
$ws.Range("F368").Value = "Utsiktens BK"
$ws.Range("G368").Value = "Orgryte IS"
$ws.Range("H368").Value = 0
$ws.Range("I368").Value = 1
$ws.Range("J368").Value = "A"
$ws.Range("K368").Value = 1.571
$ws.Range("L368").Value = 4
$ws.Range("M368").Value = 5
$ws.Range("N368").Value = 1.85
$ws.Range("O368").Value = 3.5
$ws.Range("P368").Value = 3.75
$ws.Range("Q368").Value = -0.5
$ws.Range("R368").Value = 1.95
$ws.Range("S368").Value = 1.9
$ws.Range("T368").Value = 2.75
$ws.Range("U368").Value = 2.025
$ws.Range("V368").Value = 1.825
$ws.Range("W368").Value = -1
$ws.Range("X368").Value = -1
$ws.Range("Y368").Value = 2.75
$ws.Range("Z368").Value = -1
$ws.Range("AA368").Value = 0.8999999999999999
$ws.Range("AB368").Value = -1
$ws.Range("AC368").Value = 0.825
$ws.Range("B369").Value = 4587851
$ws.Range("F369").Value = "Vasteras SK"
$ws.Range("G369").Value = "Osters IF"
$ws.Range("H369").Value = 2
$ws.Range("I369").Value = 2
$ws.Range("J369").Value = "D"
$ws.Range("K369").Value = 2.875
$ws.Range("L369").Value = 3.5
$ws.Range("M369").Value = 2.15
$ws.Range("N369").Value = 3
$ws.Range("O369").Value = 3.5
$ws.Range("P369").Value = 2.15
$ws.Range("Q369").Value = 0.25
$ws.Range("R369").Value = 1.925
$ws.Range("S369").Value = 1.925
$ws.Range("T369").Value = 2.5
$ws.Range("U369").Value = 1.825
$ws.Range("V369").Value = 2.025
$ws.Range("W369").Value = -1
$ws.Range("X369").Value = 2.5
$ws.Range("Y369").Value = -1
$ws.Range("Z369").Value = 0.4625
$ws.Range("AA369").Value = -0.5
$ws.Range("AB369").Value = 0.825
$ws.Range("AC369").Value = -1

# Rows 370, 371, 372
$ws.Range("B370").Value = 4587940
$ws.Range("F370").Value = "IK Brage"
$ws.Range("G370").Value = "Ostersunds FK"
$ws.Range("H370").Value = 1
$ws.Range("I370").Value = 0
$ws.Range("J370").Value = "H"
$ws.Range("K370").Value = 1.6
$ws.Range("L370").Value = 4
$ws.Range("M370").Value = 4.75
$ws.Range("N370").Value = 1.45
$ws.Range("O370").Value = 4.333
$ws.Range("P370").Value = 5.75
$ws.Range("Q370").Value = -1.25
$ws.Range("R370").Value = 2.05
$ws.Range("S370").Value = 1.8
$ws.Range("T370").Value = 2.75
$ws.Range("U370").Value = 1.9
$ws.Range("V370").Value = 1.95
$ws.Range("W370").Value = 0.45
$ws.Range("X370").Value = -1
$ws.Range("Y370").Value = -1
$ws.Range("Z370").Value = -0.5
$ws.Range("AA370").Value = 0.4
$ws.Range("AB370").Value = -1
$ws.Range("AC370").Value = 0.95
$ws.Range("B371").Value = 4587858
$ws.Range("F371").Value = "Trelleborgs FF"
$ws.Range("G371").Value = "Halmstad"
$ws.Range("H371").Value = 0
$ws.Range("I371").Value = 4
$ws.Range("J371").Value = "A"
$ws.Range("K371").Value = 3.25
$ws.Range("L371").Value = 3.4
$ws.Range("M371").Value = 2.05
$ws.Range("N371").Value = 3.25
$ws.Range("O371").Value = 3.4
$ws.Range("P371").Value = 2.05
$ws.Range("Q371").Value = 0.25
$ws.Range("R371").Value = 2.025
$ws.Range("S371").Value = 1.825
$ws.Range("T371").Value = 2.5
$ws.Range("U371").Value = 1.85
$ws.Range("V371").Value = 2
$ws.Range("W371").Value = -1
$ws.Range("X371").Value = -1
$ws.Range("Y371").Value = 1.05
$ws.Range("Z371").Value = -1
$ws.Range("AA371").Value = 0.825
$ws.Range("AB371").Value = 0.8500000000000001
$ws.Range("AC371").Value = -1
$ws.Range("B372").Value = 4587856
$ws.Range("F372").Value = "Landskrona BoIS"
$ws.Range("G372").Value = "Skvde AIK"
$ws.Range("H372").Value = 0
$ws.Range("I372").Value = 1
$ws.Range("J372").Value = "A"
$ws.Range("K372").Value = 2.15
$ws.Range("L372").Value = 3.4
$ws.Range("M372").Value = 3.1
$ws.Range("N372").Value = 2.55
$ws.Range("O372").Value = 3.4
$ws.Range("P372").Value = 2.55
$ws.Range("Q372").Value = 0
$ws.Range("R372").Value = 1.925
$ws.Range("S372").Value = 1.925
$ws.Range("T372").Value = 2.75
$ws.Range("U372").Value = 2.025
$ws.Range("V372").Value = 1.825
$ws.Range("W372").Value = -1
$ws.Range("X372").Value = -1
$ws.Range("Y372").Value = 1.55
$ws.Range("Z372").Value = -1
$ws.Range("AA372").Value = 0.925
$ws.Range("AB372").Value = -1
$ws.Range("AC372").Value = 0.825

# Rows 373, 374
$ws.Range("B373").Value = 4587993
$ws.Range("F373").Value = "Dalkurd FF"
$ws.Range("G373").Value = "Norrby IF"
$ws.Range("H373").Value = 2
$ws.Range("I373").Value = 1
$ws.Range("J373").Value = "H"
$ws.Range("K373").Value = 2.75
$ws.Range("L373").Value = 3.5
$ws.Range("M373").Value = 2.3
$ws.Range("N373").Value = 2.75
$ws.Range("O373").Value = 3.5
$ws.Range("P373").Value = 2.3
$ws.Range("Q373").Value = 0.25
$ws.Range("R373").Value = 1.8
$ws.Range("S373").Value = 2.05
$ws.Range("T373").Value = 2.5
$ws.Range("U373").Value = 1.825
$ws.Range("V373").Value = 2.025
$ws.Range("W373").Value = 1.75
$ws.Range("X373").Value = -1
$ws.Range("Y373").Value = -1
$ws.Range("Z373").Value = 0.8
$ws.Range("AA373").Value = -1
$ws.Range("AB373").Value = 0.825
$ws.Range("AC373").Value = -1
$ws.Range("B374").Value = 4591712
$ws.Range("F374").Value = "Brommapojkarna"
$ws.Range("G374").Value = "Orebro SK"
$ws.Range("H374").Value = 2
$ws.Range("I374").Value = 1
$ws.Range("J374").Value = "H"
$ws.Range("K374").Value = 1.85
$ws.Range("L374").Value = 3.6
$ws.Range("M374").Value = 3.6
$ws.Range("N374").Value = 1.75
$ws.Range("O374").Value = 3.75
$ws.Range("P374").Value = 3.8
$ws.Range("Q374").Value = -0.75
$ws.Range("R374").Value = 2
$ws.Range("S374").Value = 1.85
$ws.Range("T374").Value = 3
$ws.Range("U374").Value = 1.975
$ws.Range("V374").Value = 1.875
$ws.Range("W374").Value = 0.75
$ws.Range("X374").Value = -1
$ws.Range("Y374").Value = -1
$ws.Range("Z374").Value = 0.5
$ws.Range("AA374").Value = -0.5
$ws.Range("AB374").Value = 0
$ws.Range("AC374").Value = -0

# Rows 414, 415
$ws.Range("B414").Value = 4587951
$ws.Range("F414").Value = "IK Brage"
$ws.Range("G414").Value = "Utsiktens BK"
$ws.Range("H414").Value = 1
$ws.Range("I414").Value = 2
$ws.Range("J414").Value = "A"
$ws.Range("K414").Value = 2.05
$ws.Range("L414").Value = 3.3
$ws.Range("M414").Value = 3.3
$ws.Range("N414").Value = 2.2
$ws.Range("O414").Value = 3.3
$ws.Range("P414").Value = 3
$ws.Range("Q414").Value = -0.25
$ws.Range("R414").Value = 2
$ws.Range("S414").Value = 1.85
$ws.Range("T414").Value = 2.5
$ws.Range("U414").Value = 1.825
$ws.Range("V414").Value = 2.025
$ws.Range("W414").Value = -1
$ws.Range("X414").Value = -1
$ws.Range("Y414").Value = 2
$ws.Range("Z414").Value = -1
$ws.Range("AA414").Value = 0.8500000000000001
$ws.Range("AB414").Value = 0.825
$ws.Range("AC414").Value = -1
$ws.Range("B415").Value = 4587949
$ws.Range("F415").Value = "Dalkurd FF"
$ws.Range("G415").Value = "Brommapojkarna"
$ws.Range("H415").Value = 0
$ws.Range("I415").Value = 2
$ws.Range("J415").Value = "A"
$ws.Range("K415").Value = 3.6
$ws.Range("L415").Value = 3.6
$ws.Range("M415").Value = 1.833
$ws.Range("N415").Value = 3.75
$ws.Range("O415").Value = 3.75
$ws.Range("P415").Value = 1.8
$ws.Range("Q415").Value = 0.5
$ws.Range("R415").Value = 2.025
$ws.Range("S415").Value = 1.825
$ws.Range("T415").Value = 2.75
$ws.Range("U415").Value = 1.8
$ws.Range("V415").Value = 2.05
$ws.Range("W415").Value = -1
$ws.Range("X415").Value = -1
$ws.Range("Y415").Value = 0.8
$ws.Range("Z415").Value = -1
$ws.Range("AA415").Value = 0.825
$ws.Range("AB415").Value = -1
$ws.Range("AC415").Value = 1.05

# Rows 416, 417
$ws.Range("B416").Value = 4587876
$ws.Range("F416").Value = "Trelleborgs FF"
$ws.Range("G416").Value = "AFC Eskilstuna"
$ws.Range("H416").Value = 0
$ws.Range("I416").Value = 1
$ws.Range("J416").Value = "A"
$ws.Range("K416").Value = 1.7
$ws.Range("L416").Value = 3.8
$ws.Range("M416").Value = 4
$ws.Range("N416").Value = 1.533
$ws.Range("O416").Value = 4.2
$ws.Range("P416").Value = 5
$ws.Range("Q416").Value = -1
$ws.Range("R416").Value = 1.875
$ws.Range("S416").Value = 1.975
$ws.Range("T416").Value = 3.25
$ws.Range("U416").Value = 1.925
$ws.Range("V416").Value = 1.925
$ws.Range("W416").Value = -1
$ws.Range("X416").Value = -1
$ws.Range("Y416").Value = 4
$ws.Range("Z416").Value = -1
$ws.Range("AA416").Value = 0.9750000000000001
$ws.Range("AB416").Value = -1
$ws.Range("AC416").Value = 0.925
$ws.Range("B417").Value = 4587950
$ws.Range("F417").Value = "Jonkopings Sodra"
$ws.Range("G417").Value = "Osters IF"
$ws.Range("H417").Value = 2
$ws.Range("I417").Value = 1
$ws.Range("J417").Value = "H"
$ws.Range("K417").Value = 3.8
$ws.Range("L417").Value = 3.5
$ws.Range("M417").Value = 1.833
$ws.Range("N417").Value = 3.4
$ws.Range("O417").Value = 3.1
$ws.Range("P417").Value = 2.1
$ws.Range("Q417").Value = 0.25
$ws.Range("R417").Value = 2.05
$ws.Range("S417").Value = 1.8
$ws.Range("T417").Value = 2.25
$ws.Range("U417").Value = 1.875
$ws.Range("V417").Value = 1.975
$ws.Range("W417").Value = 2.4
$ws.Range("X417").Value = -1
$ws.Range("Y417").Value = -1
$ws.Range("Z417").Value = 1.05
$ws.Range("AA417").Value = -1
$ws.Range("AB417").Value = 0.875
$ws.Range("AC417").Value = -1

# Rows 424, 425
$ws.Range("B424").Value = 4587952
$ws.Range("F424").Value = "AFC Eskilstuna"
$ws.Range("G424").Value = "IK Brage"
$ws.Range("H424").Value = 2
$ws.Range("I424").Value = 3
$ws.Range("J424").Value = "A"
$ws.Range("K424").Value = 2.45
$ws.Range("L424").Value = 3.5
$ws.Range("M424").Value = 2.5
$ws.Range("N424").Value = 2.3
$ws.Range("O424").Value = 3.4
$ws.Range("P424").Value = 2.75
$ws.Range("Q424").Value = -0.25
$ws.Range("R424").Value = 2.05
$ws.Range("S424").Value = 1.75
$ws.Range("T424").Value = 2.75
$ws.Range("U424").Value = 1.85
$ws.Range("V424").Value = 2
$ws.Range("W424").Value = -1
$ws.Range("X424").Value = -1
$ws.Range("Y424").Value = 1.75
$ws.Range("Z424").Value = -1
$ws.Range("AA424").Value = 0.75
$ws.Range("AB424").Value = 0.8500000000000001
$ws.Range("AC424").Value = -1
$ws.Range("B425").Value = 4587880
$ws.Range("F425").Value = "Norrby IF"
$ws.Range("G425").Value = "Halmstad"
$ws.Range("H425").Value = 3
$ws.Range("I425").Value = 0
$ws.Range("J425").Value = "H"
$ws.Range("K425").Value = 3
$ws.Range("L425").Value = 3.5
$ws.Range("M425").Value = 2.1
$ws.Range("N425").Value = 4
$ws.Range("O425").Value = 3.6
$ws.Range("P425").Value = 1.75
$ws.Range("Q425").Value = 0.5
$ws.Range("R425").Value = 2.05
$ws.Range("S425").Value = 1.8
$ws.Range("T425").Value = 2.5
$ws.Range("U425").Value = 1.95
$ws.Range("V425").Value = 1.9
$ws.Range("W425").Value = 3
$ws.Range("X425").Value = -1
$ws.Range("Y425").Value = -1
$ws.Range("Z425").Value = 1.05
$ws.Range("AA425").Value = -1
$ws.Range("AB425").Value = 0.95
$ws.Range("AC425").Value = -1

# Rows 464, 465
$ws.Range("B464").Value = 4591701
$ws.Range("F464").Value = "Halmstad"
$ws.Range("G464").Value = "Orebro SK"
$ws.Range("H464").Value = 0
$ws.Range("I464").Value = 0
$ws.Range("J464").Value = "D"
$ws.Range("K464").Value = 1.533
$ws.Range("L464").Value = 4
$ws.Range("M464").Value = 5.25
$ws.Range("N464").Value = 1.5
$ws.Range("O464").Value = 4
$ws.Range("P464").Value = 5.75
$ws.Range("Q464").Value = -1
$ws.Range("R464").Value = 1.825
$ws.Range("S464").Value = 2.025
$ws.Range("T464").Value = 2.75
$ws.Range("U464").Value = 1.9
$ws.Range("V464").Value = 1.95
$ws.Range("W464").Value = -1
$ws.Range("X464").Value = 3
$ws.Range("Y464").Value = -1
$ws.Range("Z464").Value = -1
$ws.Range("AA464").Value = 1.025
$ws.Range("AB464").Value = -1
$ws.Range("AC464").Value = 0.95
$ws.Range("B465").Value = 4587961
$ws.Range("F465").Value = "Brommapojkarna"
$ws.Range("G465").Value = "IK Brage"
$ws.Range("H465").Value = 3
$ws.Range("I465").Value = 2
$ws.Range("J465").Value = "H"
$ws.Range("K465").Value = 1.571
$ws.Range("L465").Value = 3.9
$ws.Range("M465").Value = 5
$ws.Range("N465").Value = 1.5
$ws.Range("O465").Value = 4.2
$ws.Range("P465").Value = 5.5
$ws.Range("Q465").Value = -1
$ws.Range("R465").Value = 1.85
$ws.Range("S465").Value = 2
$ws.Range("T465").Value = 3.25
$ws.Range("U465").Value = 2
$ws.Range("V465").Value = 1.85
$ws.Range("W465").Value = 0.5
$ws.Range("X465").Value = -1
$ws.Range("Y465").Value = -1
$ws.Range("Z465").Value = 0
$ws.Range("AA465").Value = -0
$ws.Range("AB465").Value = 1
$ws.Range("AC465").Value = -1

# Rows 484, 485
$ws.Range("B484").Value = 6450192
$ws.Range("F484").Value = "Utsiktens BK"
$ws.Range("G484").Value = "Skvde AIK"
$ws.Range("H484").Value = 1
$ws.Range("I484").Value = 0
$ws.Range("J484").Value = "H"
$ws.Range("K484").Value = 2.2
$ws.Range("L484").Value = 3.3
$ws.Range("M484").Value = 3
$ws.Range("N484").Value = 1.833
$ws.Range("O484").Value = 3.4
$ws.Range("P484").Value = 4
$ws.Range("Q484").Value = -0.5
$ws.Range("R484").Value = 1.875
$ws.Range("S484").Value = 1.975
$ws.Range("T484").Value = 2.75
$ws.Range("U484").Value = 2
$ws.Range("V484").Value = 1.85
$ws.Range("W484").Value = 0.833
$ws.Range("X484").Value = -1
$ws.Range("Y484").Value = -1
$ws.Range("Z484").Value = 0.875
$ws.Range("AA484").Value = -1
$ws.Range("AB484").Value = -1
$ws.Range("AC484").Value = 0.8500000000000001
$ws.Range("B485").Value = 5992047
$ws.Range("F485").Value = "Trelleborgs FF"
$ws.Range("G485").Value = "Gefle IF"
$ws.Range("H485").Value = 0
$ws.Range("I485").Value = 1
$ws.Range("J485").Value = "A"
$ws.Range("K485").Value = 1.909
$ws.Range("L485").Value = 3.4
$ws.Range("M485").Value = 3.6
$ws.Range("N485").Value = 1.909
$ws.Range("O485").Value = 3.1
$ws.Range("P485").Value = 4
$ws.Range("Q485").Value = -0.5
$ws.Range("R485").Value = 2
$ws.Range("S485").Value = 1.85
$ws.Range("T485").Value = 2.25
$ws.Range("U485").Value = 1.95
$ws.Range("V485").Value = 1.9
$ws.Range("W485").Value = -1
$ws.Range("X485").Value = -1
$ws.Range("Y485").Value = 3
$ws.Range("Z485").Value = -1
$ws.Range("AA485").Value = 0.8500000000000001
$ws.Range("AB485").Value = -1
$ws.Range("AC485").Value = 0.8999999999999999

# Rows 584, 585
$ws.Range("B584").Value = 5992137
$ws.Range("F584").Value = "Helsingborg"
$ws.Range("G584").Value = "Ostersunds FK"
$ws.Range("H584").Value = 1
$ws.Range("I584").Value = 0
$ws.Range("J584").Value = "H"
$ws.Range("K584").Value = 2.375
$ws.Range("L584").Value = 3.2
$ws.Range("M584").Value = 2.625
$ws.Range("N584").Value = 2
$ws.Range("O584").Value = 3.4
$ws.Range("P584").Value = 3.3
$ws.Range("Q584").Value = -0.5
$ws.Range("R584").Value = 2.05
$ws.Range("S584").Value = 1.8
$ws.Range("T584").Value = 2.25
$ws.Range("U584").Value = 1.9
$ws.Range("V584").Value = 1.95
$ws.Range("W584").Value = 1
$ws.Range("X584").Value = -1
$ws.Range("Y584").Value = -1
$ws.Range("Z584").Value = 1.05
$ws.Range("AA584").Value = -1
$ws.Range("AB584").Value = -1
$ws.Range("AC584").Value = 0.95
$ws.Range("B585").Value = 5992140
$ws.Range("F585").Value = "Trelleborgs FF"
$ws.Range("G585").Value = "Utsiktens BK"
$ws.Range("H585").Value = 0
$ws.Range("I585").Value = 3
$ws.Range("J585").Value = "A"
$ws.Range("K585").Value = 2.5
$ws.Range("L585").Value = 3.4
$ws.Range("M585").Value = 2.4
$ws.Range("N585").Value = 2.75
$ws.Range("O585").Value = 3.4
$ws.Range("P585").Value = 2.25
$ws.Range("Q585").Value = 0.25
$ws.Range("R585").Value = 1.8
$ws.Range("S585").Value = 2.05
$ws.Range("T585").Value = 2.5
$ws.Range("U585").Value = 1.85
$ws.Range("V585").Value = 2
$ws.Range("W585").Value = -1
$ws.Range("X585").Value = -1
$ws.Range("Y585").Value = 1.25
$ws.Range("Z585").Value = -1
$ws.Range("AA585").Value = 1.05
$ws.Range("AB585").Value = 0.8500000000000001
$ws.Range("AC585").Value = -1

# Rows 661, 662, 663
$ws.Range("B661").Value = 5992209
$ws.Range("F661").Value = "Vasteras SK"
$ws.Range("G661").Value = "GIF Sundsvall"
$ws.Range("H661").Value = 2
$ws.Range("I661").Value = 0
$ws.Range("J661").Value = "H"
$ws.Range("K661").Value = 1.363
$ws.Range("L661").Value = 5
$ws.Range("M661").Value = 6.5
$ws.Range("N661").Value = 1.444
$ws.Range("O661").Value = 4.75
$ws.Range("P661").Value = 6.5
$ws.Range("Q661").Value = -1.25
$ws.Range("R661").Value = 1.925
$ws.Range("S661").Value = 1.925
$ws.Range("T661").Value = 3
$ws.Range("U661").Value = 1.875
$ws.Range("V661").Value = 1.975
$ws.Range("W661").Value = 0.444
$ws.Range("X661").Value = -1
$ws.Range("Y661").Value = -1
$ws.Range("Z661").Value = 0.925
$ws.Range("AA661").Value = -1
$ws.Range("AB661").Value = -1
$ws.Range("AC661").Value = 0.9750000000000001
$ws.Range("B662").Value = 5992213
$ws.Range("F662").Value = "Gefle IF"
$ws.Range("G662").Value = "GAIS"
$ws.Range("H662").Value = 0
$ws.Range("I662").Value = 6
$ws.Range("J662").Value = "A"
$ws.Range("K662").Value = 6
$ws.Range("L662").Value = 4.2
$ws.Range("M662").Value = 1.45
$ws.Range("N662").Value = 5
$ws.Range("O662").Value = 4.5
$ws.Range("P662").Value = 1.533
$ws.Range("Q662").Value = 1
$ws.Range("R662").Value = 1.925
$ws.Range("S662").Value = 1.925
$ws.Range("T662").Value = 2.75
$ws.Range("U662").Value = 1.825
$ws.Range("V662").Value = 2.025
$ws.Range("W662").Value = -1
$ws.Range("X662").Value = -1
$ws.Range("Y662").Value = 0.5329999999999999
$ws.Range("Z662").Value = -1
$ws.Range("AA662").Value = 0.925
$ws.Range("AB662").Value = 0.825
$ws.Range("AC662").Value = -1
$ws.Range("B663").Value = 5992215
$ws.Range("F663").Value = "Skvde AIK"
$ws.Range("G663").Value = "Ostersunds FK"
$ws.Range("H663").Value = 5
$ws.Range("I663").Value = 1
$ws.Range("J663").Value = "H"
$ws.Range("K663").Value = 2.875
$ws.Range("L663").Value = 3.3
$ws.Range("M663").Value = 2.25
$ws.Range("N663").Value = 3.25
$ws.Range("O663").Value = 3.4
$ws.Range("P663").Value = 2.15
$ws.Range("Q663").Value = 0.25
$ws.Range("R663").Value = 1.925
$ws.Range("S663").Value = 1.925
$ws.Range("T663").Value = 2.5
$ws.Range("U663").Value = 1.975
$ws.Range("V663").Value = 1.875
$ws.Range("W663").Value = 2.25
$ws.Range("X663").Value = -1
$ws.Range("Y663").Value = -1
$ws.Range("Z663").Value = 0.925
$ws.Range("AA663").Value = -1
$ws.Range("AB663").Value = 0.9750000000000001
$ws.Range("AC663").Value = -1

# Rows 666, 667
$ws.Range("B666").Value = 5992219
$ws.Range("F666").Value = "AFC Eskilstuna"
$ws.Range("G666").Value = "Vasteras SK"
$ws.Range("H666").Value = 1
$ws.Range("I666").Value = 2
$ws.Range("J666").Value = "A"
$ws.Range("K666").Value = 4
$ws.Range("L666").Value = 3.6
$ws.Range("M666").Value = 1.85
$ws.Range("N666").Value = 4.75
$ws.Range("O666").Value = 3.75
$ws.Range("P666").Value = 1.7
$ws.Range("Q666").Value = 0.75
$ws.Range("R666").Value = 1.875
$ws.Range("S666").Value = 1.975
$ws.Range("T666").Value = 2.75
$ws.Range("U666").Value = 1.9
$ws.Range("V666").Value = 1.95
$ws.Range("W666").Value = -1
$ws.Range("X666").Value = -1
$ws.Range("Y666").Value = 0.7
$ws.Range("Z666").Value = -0.5
$ws.Range("AA666").Value = 0.4875
$ws.Range("AB666").Value = 0.45
$ws.Range("AC666").Value = -0.5
$ws.Range("B667").Value = 5992222
$ws.Range("F667").Value = "Osters IF"
$ws.Range("G667").Value = "Helsingborg"
$ws.Range("H667").Value = 1
$ws.Range("I667").Value = 1
$ws.Range("J667").Value = "D"
$ws.Range("K667").Value = 1.444
$ws.Range("L667").Value = 4.333
$ws.Range("M667").Value = 7
$ws.Range("N667").Value = 1.444
$ws.Range("O667").Value = 4.5
$ws.Range("P667").Value = 7
$ws.Range("Q667").Value = -1.25
$ws.Range("R667").Value = 1.95
$ws.Range("S667").Value = 1.9
$ws.Range("T667").Value = 3
$ws.Range("U667").Value = 1.9
$ws.Range("V667").Value = 1.95
$ws.Range("W667").Value = -1
$ws.Range("X667").Value = 3.5
$ws.Range("Y667").Value = -1
$ws.Range("Z667").Value = -1
$ws.Range("AA667").Value = 0.8999999999999999
$ws.Range("AB667").Value = -1
$ws.Range("AC667").Value = 0.95

# Rows 679, 680
$ws.Range("B679").Value = 5993423
$ws.Range("F679").Value = "Gefle IF"
$ws.Range("G679").Value = "Orgryte IS"
$ws.Range("H679").Value = 0
$ws.Range("I679").Value = 1
$ws.Range("J679").Value = "A"
$ws.Range("K679").Value = 3.1
$ws.Range("L679").Value = 3.5
$ws.Range("M679").Value = 2.2
$ws.Range("N679").Value = 4.333
$ws.Range("O679").Value = 3.8
$ws.Range("P679").Value = 1.75
$ws.Range("Q679").Value = 0.75
$ws.Range("R679").Value = 1.875
$ws.Range("S679").Value = 1.975
$ws.Range("T679").Value = 2.75
$ws.Range("U679").Value = 1.825
$ws.Range("V679").Value = 2.025
$ws.Range("W679").Value = -1
$ws.Range("X679").Value = -1
$ws.Range("Y679").Value = 0.75
$ws.Range("Z679").Value = -0.5
$ws.Range("AA679").Value = 0.4875
$ws.Range("AB679").Value = -1
$ws.Range("AC679").Value = 1.025
$ws.Range("B680").Value = 5992227
$ws.Range("F680").Value = "Landskrona BoIS"
$ws.Range("G680").Value = "GIF Sundsvall"
$ws.Range("H680").Value = 2
$ws.Range("I680").Value = 2
$ws.Range("J680").Value = "D"
$ws.Range("K680").Value = 1.909
$ws.Range("L680").Value = 3.75
$ws.Range("M680").Value = 3.8
$ws.Range("N680").Value = 1.85
$ws.Range("O680").Value = 3.8
$ws.Range("P680").Value = 3.75
$ws.Range("Q680").Value = -0.5
$ws.Range("R680").Value = 1.9
$ws.Range("S680").Value = 1.95
$ws.Range("T680").Value = 3
$ws.Range("U680").Value = 1.925
$ws.Range("V680").Value = 1.925
$ws.Range("W680").Value = -1
$ws.Range("X680").Value = 2.8
$ws.Range("Y680").Value = -1
$ws.Range("Z680").Value = -1
$ws.Range("AA680").Value = 0.95
$ws.Range("AB680").Value = 0.925
$ws.Range("AC680").Value = -1
